$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Dong%Liu%NULL%1,                    Roy F.%Baumeister%NULL%1,                    Jennifer C.%Veilleux%NULL%1,                    Caixia%Chen%NULL%1,                    Wenjun%Liu%NULL%1,                    Yongjie%Yue%NULL%1,                    Shi%Zhang%NULL%1]'
$ws.Range("I2").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("I3").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E4").Value = '[Xin%Cai%NULL%1,                    Xiaopeng%Hu%NULL%1,                    Ivo Otte%Ekumi%NULL%1,                    Jianchun%Wang%NULL%1,                    Yawen%An%NULL%1,                    Zhiwen%Li%NULL%1,                    Bo%Yuan%NULL%1]'
$ws.Range("I4").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E5").Value = '[Hai-Xin%Bo%NULL%0,                    Wen%Li%NULL%2,                    Yuan%Yang%NULL%4,                    Yuan%Yang%NULL%0,                    Yu%Wang%NULL%2,                    Qinge%Zhang%NULL%1,                    Teris%Cheung%NULL%2,                    Xinjuan%Wu%NULL%2,                    Yu-Tao%Xiang%NULL%2]'
$ws.Range("I5").Value = '_PMC_CROSSREF'
$ws.Range("E6").Value = '[Michael L.%Tee%NULL%1,                    Cherica A.%Tee%NULL%1,                    Joseph P.%Anlacan%NULL%1,                    Katrina Joy G.%Aligam%NULL%1,                    Patrick Wincy C.%Reyes%NULL%1,                    Vipat%Kuruchittham%NULL%1,                    Roger C.%Ho%NULL%2]'
$ws.Range("I6").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E7").Value = '[Ming-Yu%Si%myu_si@163.com%1,                    Xiao-You%Su%suxiaoyou@hotmail.com%1,                    Yu%Jiang%jiangyu@pumc.edu.cn%2,                    Yu%Jiang%jiangyu@pumc.edu.cn%0,                    Wen-Jun%Wang%wwjun1973@163.com%1,                    Xiao-Fen%Gu%357539328@qq.com%1,                    Li%Ma%mali_lele@sina.com%1,                    Jing%Li%lijingwcph@scu.edu.cn%0,                    Shao-Kai%Zhang%shaokaizhang@126.com%1,                    Ze-Fang%Ren%renzef@mail.sysu.edu.cn%1,                    Ran%Ren%renran99@163.com%1,                    Yuan-Li%Liu%liuyuanli_pumc@163.com%1,                    You-Lin%Qiao%qiaoy@cicams.ac.cn%1]'
$ws.Range("I7").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E8").Value = '[Rocío%Rodríguez-Rey%NULL%1,                    Helena%Garrido-Hernansaiz%NULL%1,                    Silvia%Collado%NULL%1]'
$ws.Range("I8").Value = '_PMC_CROSSREF'
$ws.Range("E9").Value = '[Anliu%Nie%NULL%1,                    Xiangfen%Su%suxiangfengfy@163.com%1,                    Shuzeng%Zhang%NULL%2,                    Shuzeng%Zhang%NULL%0,                    Wenjie%Guan%NULL%1,                    Jianfeng%Li%NULL%1]'
$ws.Range("I9").Value = '_PMC_CROSSREF'
$ws.Range("E10").Value = '[Shun-Wei%Liang%NULL%1,                    Rong-Ning%Chen%NULL%1,                    Li-Li%Liu%NULL%1,                    Xue-Guo%Li%NULL%1,                    Jian-Bin%Chen%NULL%1,                    Si-Yao%Tang%NULL%1,                    Jing-Bo%Zhao%NULL%1]'
$ws.Range("I10").Value = '_PMC_CROSSREF'
$ws.Range("E11").Value = '[Guo%Li%xref no email%2,            Jinfeng%Miao%xref no email%2,            Hui%Wang%xref no email%2,            Shabei%Xu%xref no email%2,            Wenzhe%Sun%xref no email%2,            Yebin%Fan%xref no email%2,            Chenyan%Zhang%xref no email%2,            Suiqiang%Zhu%xref no email%2,            Zhou%Zhu%xref no email%2,            Wei%Wang%xref no email%0]'
$ws.Range("I11").Value = '_CROSSREF'
$ws.Range("E12").Value = '[Emanuele Maria%Giusti%NULL%2,                    Elisa%Pedroli%NULL%1,                    Guido E.%D''Aniello%NULL%1,                    Chiara%Stramba Badiale%NULL%1,                    Giada%Pietrabissa%NULL%1,                    Chiara%Manna%NULL%1,                    Marco%Stramba Badiale%NULL%1,                    Giuseppe%Riva%NULL%1,                    Gianluca%Castelnuovo%NULL%1,                    Enrico%Molinari%NULL%1]'
$ws.Range("I12").Value = '_PMC_CROSSREF'
$ws.Range("E13").Value = '[Biao%Chen%NULL%1,                    Qing-xian%Li%NULL%1,                    Heng%Zhang%NULL%1,                    Jia-yong%Zhu%NULL%1,                    Xu%Yang%NULL%1,                    Yu-hang%Wu%NULL%1,                    Jie%Xiong%NULL%1,                    Fu%Li%NULL%1,                    Hua%Wang%whzkwlh@126.com%0,                    Zhi-tao%Chen%chenzhitao830201@163.com%1]'
$ws.Range("I13").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E14").Value = '[Anaëlle%Caillet%NULL%1,                    Charlotte%Coste%NULL%1,                    Rocio%Sanchez%NULL%1,                    Bernard%Allaouchiche%NULL%1]'
$ws.Range("I14").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E15").Value = '[Mariapaola%Barbato%Mariapaola.Barbato@zu.ac.ae%1,                    Justin%Thomas%NULL%1]'
$ws.Range("I15").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E16").Value = '[Abdulmajeed A.%Alkhamees%NULL%1,                    Saleh A.%Alrashed%NULL%1,                    Ali A.%Alzunaydi%NULL%1,                    Ahmed S.%Almohimeed%NULL%1,                    Moath S.%Aljohani%NULL%1]'
$ws.Range("I16").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E17").Value = '[Yongjie%Zhou%NULL%1,                    Hui%Shi%NULL%1,                    Zhengkui%Liu%NULL%1,                    Songxu%Peng%gwxypsx@163.com%2,                    Songxu%Peng%gwxypsx@163.com%0,                    Ruoxi%Wang%NULL%2,                    Ruoxi%Wang%NULL%0,                    Ling%Qi%NULL%1,                    Zezhi%Li%NULL%1,                    Jiezhi%Yang%NULL%1,                    Yali%Ren%NULL%1,                    Xiuli%Song%NULL%1,                    Lingyun%Zeng%NULL%1,                    Wei%Qian%NULL%1,                    Xiangyang%Zhang%zhangxy@psych.ac.cn%1]'
$ws.Range("I17").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E18").Value = '[Chang%J.%coreGivesNoEmail%1,                   Khazaie%H.%coreGivesNoEmail%1,                   Li%J.%coreGivesNoEmail%1,                   Xiaohui%Li%coreGivesNoEmail%1,                   Xing%Tan%coreGivesNoEmail%1,                   Yuanyuan%An%coreGivesNoEmail%1,                   Yuqing%Zhao%coreGivesNoEmail%1]'
$ws.Range("I18").Value = '_CROSSREF'
$ws.Range("E19").Value = '[Yingfei%Zhang%NULL%0,                    Zheng Feei%Ma%NULL%4,                    Zheng Feei%Ma%NULL%0]'
$ws.Range("I19").Value = '_PMC_CROSSREF'
$ws.Range("E20").Value = '[Qianlan%Yin%NULL%1,                    Zhuoer%Sun%NULL%3,                    Zhuoer%Sun%NULL%0,                    Tuanjie%Liu%NULL%1,                    Xiong%Ni%NULL%1,                    Xuanfeng%Deng%NULL%1,                    Yanpu%Jia%NULL%2,                    Zhilei%Shang%NULL%2,                    Yaoguang%Zhou%NULL%2,                    Weizhi%Liu%13024141970@163.com%3,                    Weizhi%Liu%13024141970@163.com%0]'
$ws.Range("I20").Value = '_PMC_CROSSREF'
$ws.Range("E21").Value = '[U.%Wesemann%NULL%1,                    N.%Hadjamu%NULL%2,                    N.%Hadjamu%NULL%0,                    G.%Willmund%NULL%1,                    S.%Dolff%NULL%1,                    N.%Vonderlin%NULL%1,                    R.%Wakili%NULL%1,                    J.%Vogel%NULL%1,                    T.%Rassaf%NULL%1,                    J.%Siebermair%NULL%1]'
$ws.Range("I21").Value = '_PMC_CROSSREF'
$ws.Range("E22").Value = '[Cuiyan%Wang%NULL%0,                    Riyu%Pan%NULL%1,                    Xiaoyang%Wan%NULL%1,                    Yilin%Tan%NULL%1,                    Linkang%Xu%NULL%1,                    Roger S.%McIntyre%NULL%3,                    Faith N.%Choo%NULL%1,                    Bach%Tran%NULL%3,                    Roger%Ho%NULL%3,                    Vijay K.%Sharma%NULL%2,                    Cyrus%Ho%NULL%3]'
$ws.Range("I22").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E23").Value = '[Mohit%Varshney%NULL%1,                    Jithin Thomas%Parel%NULL%2,                    Jithin Thomas%Parel%NULL%0,                    Neeraj%Raizada%NULL%2,                    Neeraj%Raizada%NULL%0,                    Shiv Kumar%Sarin%NULL%1,                    Abdallah M.%Samy%NULL%2,                    Abdallah M.%Samy%NULL%0]'
$ws.Range("I23").Value = '_PMC_CROSSREF'
$ws.Range("E24").Value = '[Claudia%Traunmüller%claudia.traunmueller@uni-graz.at%1,                    Rene%Stefitz%NULL%1,                    Kerstin%Gaisbachgrabner%NULL%1,                    Andreas%Schwerdtfeger%NULL%1]'
$ws.Range("I24").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E25").Value = '[Wanjie%Tang%NULL%1,                    Tao%Hu%NULL%1,                    Baodi%Hu%NULL%1,                    Chunhan%Jin%NULL%1,                    Gang%Wang%NULL%1,                    Chao%Xie%NULL%1,                    Sen%Chen%NULL%1,                    Jiuping%Xu%NULL%1]'
$ws.Range("I25").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E26").Value = '[Wanqiu%Tan%NULL%1,                    Fengyi%Hao%NULL%1,                    Roger S.%McIntyre%NULL%0,                    Li%Jiang%NULL%2,                    Xiaojiang%Jiang%NULL%2,                    Ling%Zhang%NULL%2,                    Xinling%Zhao%NULL%2,                    Yiran%Zou%NULL%2,                    Yirong%Hu%NULL%2,                    Xi%Luo%NULL%2,                    Zhisong%Zhang%NULL%2,                    Andre%Lai%NULL%1,                    Roger%Ho%NULL%0,                    Bach%Tran%NULL%0,                    Cyrus%Ho%NULL%0,                    Wilson%Tam%NULL%2]'
$ws.Range("I26").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E27").Value = '[Xingyue%Song%NULL%3,                    Wenning%Fu%NULL%2,                    Xiaoran%Liu%NULL%2,                    Zhiqian%Luo%NULL%2,                    Rixing%Wang%NULL%2,                    Ning%Zhou%NULL%2,                    Shijiao%Yan%NULL%3,                    Chuanzhu%Lv%NULL%3]'
$ws.Range("I27").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E28").Value = '[Allen C.%Sherman%NULL%1,                    Mark L%Williams%NULL%1,                    Benjamin C%Amick%NULL%1,                    Teresa J.%Hudson%NULL%1,                    Erick L%Messias%NULL%1]'
$ws.Range("I28").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E29").Value = '[Emire%Seyahi%eseyahi@yahoo.com%1,                    Burc Cagri%Poyraz%bcpoyraz@yahoo.com%2,                    Burc Cagri%Poyraz%bcpoyraz@yahoo.com%0,                    Necdet%Sut%necdetsut@yahoo.com%2,                    Necdet%Sut%necdetsut@yahoo.com%0,                    Selma%Akdogan%drslma.akdogan@gmail.com%2,                    Selma%Akdogan%drslma.akdogan@gmail.com%0,                    Vedat%Hamuryudan%vhamuryudan@yahoo.com%2,                    Vedat%Hamuryudan%vhamuryudan@yahoo.com%0]'
$ws.Range("I29").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E30").Value = '[Rodolfo%Rossi%NULL%3,                    Valentina%Socci%NULL%3,                    Dalila%Talevi%NULL%1,                    Sonia%Mensi%NULL%1,                    Cinzia%Niolu%NULL%1,                    Francesca%Pacitti%NULL%3,                    Antinisca%Di Marco%NULL%3,                    Alessandro%Rossi%NULL%3,                    Alberto%Siracusano%NULL%3,                    Giorgio%Di Lorenzo%NULL%3]'
$ws.Range("I30").Value = '_PMC_CROSSREF'
$ws.Range("E31").Value = '[Rodolfo%Rossi%NULL%0,                    Valentina%Socci%NULL%0,                    Francesca%Pacitti%NULL%0,                    Giorgio%Di Lorenzo%NULL%0,                    Antinisca%Di Marco%NULL%0,                    Alberto%Siracusano%NULL%0,                    Alessandro%Rossi%NULL%0]'
$ws.Range("I31").Value = '_PMC_CROSSREF'
$ws.Range("E32").Value = '[Marianna%Riello%NULL%1,                    Marianna%Purgato%NULL%1,                    Chiara%Bove%NULL%1,                    David%MacTaggart%NULL%1,                    Elena%Rusconi%NULL%1]'
$ws.Range("I32").Value = '_PMC_CROSSREF'
$ws.Range("E33").Value = '[Rongfeng%Qi%NULL%2,                    Wei%Chen%NULL%2,                    Saiduo%Liu%NULL%2,                    Paul M.%Thompson%NULL%2,                    Long Jiang%Zhang%NULL%2,                    Fei%Xia%NULL%3,                    Fang%Cheng%NULL%3,                    Ailing%Hong%NULL%2,                    Wesley%Surento%NULL%2,                    Song%Luo%NULL%2,                    Zhi Yuan%Sun%NULL%2,                    Chang Sheng%Zhou%NULL%2,                    Lingjiang%Li%NULL%2,                    Xiangao%Jiang%NULL%2,                    Guang Ming%Lu%NULL%2,                     R.%Qi%null%2,                     W.% Chen%null%2,                     S.% Liu%null%2,                     P. M.% Thompson%null%2,                     L. J.% Zhang%null%2,                     F.% Xia%null%2,                     F.% Cheng%null%2,                     A.% Hong%null%2,                     W.% Surento%null%2,                     S.% Luo%null%2,                     Z. Y.% Sun%null%2,                     C. S.% Zhou%null%2,                     L.% Li%null%2,                     X.% Jiang%null%2,                     G. M. % Lu%null%2,                   R.%Qi%null%1,                   W.% Chen%null%1,                   S.% Liu%null%1,                   P. M.% Thompson%null%1,                   L. J.% Zhang%null%1,                   F.% Xia%null%1,                   F.% Cheng%null%1,                   A.% Hong%null%1,                   W.% Surento%null%1,                   S.% Luo%null%1,                   Z. Y.% Sun%null%1,                   C. S.% Zhou%null%1,                   L.% Li%null%1,                   X.% Jiang%null%1,                   G. M. % Lu%null%1]'
$ws.Range("I33").Value = '_PMC_CROSSREF'
$ws.Range("E34").Value = '[Zheng Feei%Ma%NULL%1,                    Yutong%Zhang%NULL%2,                    Yutong%Zhang%NULL%0,                    Xiaoqin%Luo%NULL%1,                    Xinli%Li%NULL%1,                    Yeshan%Li%NULL%1,                    Shuchang%Liu%NULL%1,                    Yingfei%Zhang%NULL%1]'
$ws.Range("I34").Value = '_PMC_CROSSREF'
$ws.Range("E35").Value = '[Lourdes%Luceño-Moreno%NULL%1,                    Beatriz%Talavera-Velasco%NULL%2,                    Beatriz%Talavera-Velasco%NULL%0,                    Yolanda%García-Albuerne%NULL%1,                    Jesús%Martín-García%NULL%1]'
$ws.Range("I35").Value = '_PMC_CROSSREF'
$ws.Range("E36").Value = '[Nianqi%Liu%NULL%0,                    Fan%Zhang%NULL%1,                    Cun%Wei%NULL%1,                    Yanpu%Jia%NULL%0,                    Zhilei%Shang%NULL%0,                    Luna%Sun%NULL%1,                    Lili%Wu%NULL%1,                    Zhuoer%Sun%NULL%0,                    Yaoguang%Zhou%NULL%0,                    Yan%Wang%NULL%1,                    Weizhi%Liu%13024141970@163.com%0]'
$ws.Range("I36").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E37").Value = '[Cindy H.%Liu%chliu@bwh.harvard.edu%1,                    Emily%Zhang%NULL%1,                    Ga Tin Fifi%Wong%NULL%1,                    Sunah%Hyun%NULL%1,                    Hyeouk “Chris”%Hahm%NULL%1]'
$ws.Range("I37").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E38").Value = '[Yuchen%Li%NULL%1,                    Yue%Wang%NULL%1,                    Jingwen%Jiang%NULL%1,                    Unnur A.%Valdimarsdóttir%NULL%1,                    Katja%Fall%NULL%2,                    Katja%Fall%NULL%0,                    Fang%Fang%NULL%4,                    Fang%Fang%NULL%0,                    Huan%Song%NULL%1,                    Donghao%Lu%NULL%1,                    Wei%Zhang%NULL%2]'
$ws.Range("I38").Value = '_PMC_CROSSREF'
$ws.Range("I39").Value = '_CROSSREF'
$ws.Range("E40").Value = '[Bisson%J. I.%coreGivesNoEmail%1,                   Chiang%E. S.%coreGivesNoEmail%1,                   Fernandez%R.%coreGivesNoEmail%1,                   Gonz\u00e1lez%V. B. A.%coreGivesNoEmail%1,                   Ibrahim%H.%coreGivesNoEmail%1,                   Lee%B.%coreGivesNoEmail%1,                   Shi%L.%coreGivesNoEmail%1,                   Sun%N.%coreGivesNoEmail%1,                   Tsay%S.\u2010F.%coreGivesNoEmail%1,                   World%Health Organisation%coreGivesNoEmail%2,                   World%Health Organisation%coreGivesNoEmail%0,                   Xiaoyue%Xu%coreGivesNoEmail%1,                   Xiuchuan%Li%coreGivesNoEmail%1,                   Ying%Zhou%coreGivesNoEmail%3]'
$ws.Range("I40").Value = '_CROSSREF'
$ws.Range("E41").Value = '[Xin%Li%NULL%1,                    Peixin%Lu%NULL%1,                    Lianting%Hu%NULL%1,                    Tianhui%Huang%NULL%1,                    Long%Lu%NULL%1]'
$ws.Range("I41").Value = '_PMC_CROSSREF'
$ws.Range("I42").Value = '_PMC_CROSSREF'
$ws.Range("E43").Value = '[Guorong%Cao%coreGivesNoEmail%1,                   Holly%Wei%coreGivesNoEmail%1,                   Hong%Xu%coreGivesNoEmail%1,                   Lili%Wei%coreGivesNoEmail%1,                   Lu%H%coreGivesNoEmail%1,                   Min%Leng%coreGivesNoEmail%1,                   Shi%Ti JC%coreGivesNoEmail%1,                   Shuyun%Xing%coreGivesNoEmail%1,                   Weathers%F%coreGivesNoEmail%1,                   Wenwen%Zhang%coreGivesNoEmail%1,                   Xiaohui%Shi%coreGivesNoEmail%1,                   Xiaoying%Zhang%coreGivesNoEmail%1,                   Yang%T%coreGivesNoEmail%1,                   Yuling%Wei%coreGivesNoEmail%1]'
$ws.Range("I43").Value = '_CROSSREF'
$ws.Range("E44").Value = '[Xuan Thi Thanh%Le%NULL%1,                    Anh Kim%Dang%NULL%1,                    Jayson%Toweh%NULL%1,                    Quang Nhat%Nguyen%NULL%1,                    Huong Thi%Le%NULL%1,                    Toan Thi Thanh%Do%NULL%1,                    Hanh Bich Thi%Phan%NULL%1,                    Thao Thanh%Nguyen%NULL%1,                    Quan Thi%Pham%NULL%1,                    Nhung Kim Thi%Ta%NULL%1,                    Quynh Thi%Nguyen%NULL%1,                    Anh Ngoc%Nguyen%NULL%1,                    Quan%Van Duong%NULL%1,                    Men Thi%Hoang%NULL%1,                    Hai Quang%Pham%NULL%1,                    Linh Gia%Vu%NULL%1,                    Bach Xuan%Tran%NULL%1,                    Carl A.%Latkin%NULL%1,                    Cyrus S. H.%Ho%NULL%1,                    Roger C. M.%Ho%NULL%1]'
$ws.Range("I44").Value = '_PMC_CROSSREF'
$ws.Range("E45").Value = '[M.%Lange%NULL%1,                    S.%Joo%NULL%1,                    P.-A.%Couette%NULL%1,                    S.%de Jaegher%NULL%1,                    F.%Joly%NULL%1,                    X.%Humbert%NULL%1]'
$ws.Range("I45").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E46").Value = '[Jianbo%Lai%NULL%0,                    Simeng%Ma%NULL%2,                    Ying%Wang%NULL%4,                    Zhongxiang%Cai%NULL%2,                    Jianbo%Hu%NULL%2,                    Ning%Wei%NULL%2,                    Jiang%Wu%NULL%3,                    Hui%Du%NULL%2,                    Tingting%Chen%NULL%2,                    Ruiting%Li%NULL%2,                    Huawei%Tan%NULL%2,                    Lijun%Kang%NULL%2,                    Lihua%Yao%NULL%2,                    Manli%Huang%NULL%2,                    Huafen%Wang%NULL%2,                    Gaohua%Wang%NULL%3,                    Zhongchun%Liu%NULL%3,                    Shaohua%Hu%NULL%2]'
$ws.Range("I46").Value = '_PMC_CROSSREF'
$ws.Range("E47").Value = '[Thanos%Karatzias%t.karatzias@napier.ac.uk%1,                    Mark%Shevlin%NULL%1,                    Jamie%Murphy%NULL%1,                    Orla%McBride%NULL%1,                    Menachem%Ben‐Ezra%NULL%2,                    Menachem%Ben‐Ezra%NULL%0,                    Richard P.%Bentall%NULL%1,                    Frédérique%Vallières%NULL%1,                    Philip%Hyland%NULL%1]'
$ws.Range("I47").Value = '_PMC_CROSSREF'
$ws.Range("E48").Value = '[Michelle I.%Cardel%mcardel@ufl.edu%1,                    Stephanie%Manasse%NULL%2,                    Stephanie%Manasse%NULL%0,                    Rebecca A.%Krukowski%NULL%2,                    Rebecca A.%Krukowski%NULL%0,                    Kathryn%Ross%NULL%2,                    Kathryn%Ross%NULL%0,                    Rebecca%Shakour%NULL%2,                    Rebecca%Shakour%NULL%0,                    Darci R.%Miller%NULL%1,                    Dominick J.%Lemas%NULL%1,                    Young‐Rock%Hong%NULL%1]'
$ws.Range("I48").Value = '_PMC_CROSSREF'
$ws.Range("E49").Value = '[Jing%Guo%NULL%1,                    Xing Lin%Feng%NULL%1,                    Xiao Hua%Wang%NULL%2,                    Xiao Hua%Wang%NULL%0,                    Marinus H.%van IJzendoorn%NULL%1]'
$ws.Range("I49").Value = '_PMC_CROSSREF'
$ws.Range("E50").Value = '[Clara%González-Sanguino%NULL%2,                    Berta%Ausín%NULL%1,                    Miguel Ángel%Castellanos%NULL%1,                    Jesús%Saiz%NULL%1,                    Aída%López-Gómez%NULL%1,                    Carolina%Ugidos%NULL%1,                    Manuel%Muñoz%NULL%1]'
$ws.Range("I50").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E51").Value = '[Leivy Patricia%González Ramírez%NULL%1,                    Reyna Jazmín%Martínez Arriaga%NULL%2,                    Reyna Jazmín%Martínez Arriaga%NULL%0,                    Martha Alicia%Hernández-Gonzalez%NULL%2,                    Martha Alicia%Hernández-Gonzalez%NULL%0,                    José María%De la Roca-Chiapas%NULL%2,                    José María%De la Roca-Chiapas%NULL%0]'
$ws.Range("I51").Value = '_PMC_CROSSREF'
$ws.Range("E52").Value = '[Giuseppe%Forte%NULL%1,                    Francesca%Favieri%NULL%2,                    Francesca%Favieri%NULL%0,                    Renata%Tambelli%NULL%2,                    Renata%Tambelli%NULL%0,                    Maria%Casagrande%NULL%1]'
$ws.Range("I52").Value = '_PMC_CROSSREF'
$ws.Range("E53").Value = '[Feten%Fekih-Romdhane%feten.fekih@gmail.com%1,                    Farah%Ghrissi%NULL%1,                    Bouthaina%Abbassi%NULL%1,                    Wissal%Cherif%NULL%1,                    Majda%Cheour%NULL%1]'
$ws.Range("I53").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E54").Value = '[Safaa M.%El-Zoghby%safaaelzoghby86@hotmail.com%1,                    Enayat M.%Soltan%enayatsoltan@gmail.com%2,                    Enayat M.%Soltan%enayatsoltan@gmail.com%0,                    Hend M.%Salama%hind_mikhail@yahoo.com%2,                    Hend M.%Salama%hind_mikhail@yahoo.com%0]'
$ws.Range("I54").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E55").Value = '[Hannah%Dobson%NULL%1,                    Charles B%Malpas%NULL%2,                    Charles B%Malpas%NULL%0,                    Aidan JC%Burrell%NULL%1,                    Caroline%Gurvich%NULL%1,                    Leo%Chen%NULL%1,                    Jayashri%Kulkarni%NULL%1,                    Toby%Winton-Brown%NULL%2,                    Toby%Winton-Brown%NULL%0]'
$ws.Range("I55").Value = '_PMC_CROSSREF'
$ws.Range("E56").Value = '[Marialaura%Di Tella%xref no email%1,            Annunziata%Romeo%xref no email%1,            Agata%Benfante%xref no email%1,            Lorys%Castelli%xref no email%1]'
$ws.Range("I56").Value = '_CROSSREF'
$ws.Range("E57").Value = '[Nadia Yanet%Cortés-Álvarez%NULL%1,                    Regino%Piñeiro-Lamas%NULL%1,                    César Rubén%Vuelvas-Olmos%NULL%1]'
$ws.Range("I57").Value = '_PMC_CROSSREF'
$ws.Range("E58").Value = '[Alyssa M.%Civantos%NULL%2,                    Yasmeen%Byrnes%NULL%2,                    Yasmeen%Byrnes%NULL%0,                    Changgee%Chang%NULL%2,                    Aman%Prasad%NULL%1,                    Kevin%Chorath%NULL%2,                    Kevin%Chorath%NULL%0,                    Seerat K.%Poonia%NULL%1,                    Carolyn M.%Jenks%NULL%1,                    Andrés M.%Bur%NULL%1,                    Punam%Thakkar%NULL%2,                    Punam%Thakkar%NULL%0,                    Evan M.%Graboyes%NULL%2,                    Evan M.%Graboyes%NULL%0,                    Rahul%Seth%NULL%2,                    Rahul%Seth%NULL%0,                    Samuel%Trosman%NULL%1,                    Anni%Wong%NULL%2,                    Anni%Wong%NULL%0,                    Benjamin M.%Laitman%NULL%1,                    Brianna N.%Harris%NULL%1,                    Janki%Shah%NULL%2,                    Janki%Shah%NULL%0,                    Vanessa%Stubbs%NULL%1,                    Garret%Choby%NULL%1,                    Qi%Long%NULL%2,                    Christopher H.%Rassekh%NULL%1,                    Erica%Thaler%NULL%1,                    Karthik%Rajasekaran%karthik.rajasekaran@pennmedicine.upenn.edu%2]'
$ws.Range("I58").Value = '_PMC_CROSSREF'
$ws.Range("E59").Value = '[Alyssa M.%Civantos%NULL%0,                    Antonio%Bertelli%NULL%1,                    Antonio%Gonçalves%NULL%1,                    Emily%Getzen%NULL%1,                    Changgee%Chang%NULL%0,                    Qi%Long%NULL%0,                    Karthik%Rajasekaran%NULL%0]'
$ws.Range("I59").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E60").Value = '[Xinli%Chi%NULL%2,                    Benjamin%Becker%NULL%2,                    Qian%Yu%NULL%2,                    Peter%Willeit%NULL%2,                    Can%Jiao%NULL%2,                    Liuyue%Huang%NULL%2,                    M. Mahhub%Hossain%NULL%2,                    Igor%Grabovac%NULL%2,                    Albert%Yeung%NULL%2,                    Jingyuan%Lin%NULL%2,                    Nicola%Veronese%NULL%2,                    Jian%Wang%NULL%3,                    Xinqi%Zhou%NULL%2,                    Scott R.%Doig%NULL%2,                    Xiaofeng%Liu%NULL%2,                    Andre F.%Carvalho%NULL%2,                    Lin%Yang%NULL%4,                    Tao%Xiao%NULL%2,                    Liye%Zou%NULL%2,                    Paolo%Fusar-Poli%NULL%2,                    Marco%Solmi%NULL%2]'
$ws.Range("I60").Value = '_PMC_CROSSREF'
$ws.Range("E61").Value = '[Nicholas W. S.%Chew%NULL%1,                    Jinghao Nicholas%Ngiam%NULL%1,                    Benjamin Yong-Qiang%Tan%NULL%2,                    Benjamin Yong-Qiang%Tan%NULL%0,                    Sai-Meng%Tham%NULL%1,                    Celine Yan-Shan%Tan%NULL%1,                    Mingxue%Jing%NULL%1,                    Renarebecca%Sagayanathan%NULL%1,                    Jin Tao%Chen%NULL%1,                    Lily Y. H.%Wong%NULL%1,                    Aftab%Ahmad%NULL%1,                    Faheem Ahmed%Khan%NULL%1,                    Maznah%Marmin%NULL%1,                    Fadhlina Binte%Hassan%NULL%1,                    Tai Mei-Ling%Sharon%NULL%1,                    Chin Han%Lim%NULL%1,                    Mohamad Iqbal Bin%Mohaini%NULL%1,                    Rivan%Danuaji%NULL%1,                    Thang H.%Nguyen%NULL%1,                    Georgios%Tsivgoulis%NULL%1,                    Sotirios%Tsiodras%NULL%1,                    Paraskevi C.%Fragkou%NULL%1,                    Dimitra%Dimopoulou%NULL%2,                    Dimitra%Dimopoulou%NULL%0,                    Arvind K.%Sharma%NULL%1,                    Kenam%Shah%NULL%1,                    Bhargesh%Patel%NULL%1,                    Suktara%Sharma%NULL%1,                    R. N.%Komalkumar%NULL%1,                    R. V.%Meenakshi%NULL%1,                    Shikha%Talati%NULL%1,                    Hock Luen%Teoh%NULL%1,                    Cyrus S.%Ho%NULL%1,                    Roger C.%Ho%NULL%0,                    Vijay K.%Sharma%NULL%0]'
$ws.Range("I61").Value = '_PMC_CROSSREF'
$ws.Range("E62").Value = '[Zhongxiang%Cai%NULL%1,                    Qin%Cui%NULL%1,                    Zhongchun%Liu%NULL%0,                    Juanjuan%Li%NULL%1,                    Xuan%Gong%NULL%1,                    Jingfang%Liu%NULL%1,                    Zhiying%Wan%NULL%1,                    Xiaoping%Yuan%NULL%1,                    Xiaofen%Li%NULL%1,                    Chuang%Chen%NULL%1,                    Gaohua%Wang%NULL%0]'
$ws.Range("I62").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E63").Value = '[Aliki%F. Chatzigeorgiou%coreGivesNoEmail%1,                   Apostolos%Blekas%coreGivesNoEmail%1,                   Eleni%Parlapani%coreGivesNoEmail%1,                   Ioannis%Diakogiannis%coreGivesNoEmail%1,                   Maria%Athanasiadou%coreGivesNoEmail%1,                   Marina%Skoupra%coreGivesNoEmail%1,                   Markos%Syngelakis%coreGivesNoEmail%1,                   Panteleimon%Voitsidis%coreGivesNoEmail%1,                   Vasiliki%Holeva%coreGivesNoEmail%2]'
$ws.Range("I63").Value = '_CROSSREF'
$ws.Range("E64").Value = '[Yingfei%Zhang%NULL%0,                    Zheng Feei%Ma%NULL%0,                    Zheng Feei%Ma%NULL%0]'
$ws.Range("I64").Value = '_PMC_CROSSREF'
$ws.Range("E65").Value = '[Fengyi%Hao%NULL%1,                    Wanqiu%Tan%NULL%1,                    Li%Jiang%NULL%0,                    Ling%Zhang%NULL%0,                    Xinling%Zhao%NULL%0,                    Yiran%Zou%NULL%0,                    Yirong%Hu%NULL%0,                    Xi%Luo%NULL%0,                    Xiaojiang%Jiang%NULL%0,                    Roger S.%McIntyre%NULL%0,                    Bach%Tran%NULL%0,                    Jiaqian%Sun%NULL%1,                    Zhisong%Zhang%NULL%0,                    Roger%Ho%NULL%0,                    Cyrus%Ho%NULL%0,                    Wilson%Tam%NULL%0]'
$ws.Range("I65").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E66").Value = '[Leilei%Liang%NULL%1,                    Hui%Ren%NULL%1,                    Ruilin%Cao%NULL%1,                    Yueyang%Hu%NULL%1,                    Zeying%Qin%NULL%1,                    Chuanen%Li%NULL%1,                    Songli%Mei%meisongli@sina.com%1]'
$ws.Range("I66").Value = '_PMC_Springer_CROSSREF'
